$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new date header in BB1, copying the date/number formatting (style)
# from BA1 so the new column matches the existing header look (border, bold,
# centered, custom date format).
$ws.Range("BA1").Copy() | Out-Null
$ws.Range("BB1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("BB1").Value2 = 45986

# Column BB values for rows 3-18 simply repeat the value already present in
# column BA for that row.
$repeatRows = 3..18
foreach ($r in $repeatRows) {
    $baVal = $ws.Range("BA$r").Value2
    if ($baVal -ne $null) {
        $ws.Range("BB$r").Value2 = $baVal
    }
}

# Rows 19-21 get new, distinct forecast values in column BB.
$ws.Range("BB19").Value2 = 1.049317648994741
$ws.Range("BB20").Value2 = 0.72625340902297
$ws.Range("BB21").Value2 = 0.8024032016000104
